# GHB instruments and deployments
# Adds two new instrument models (Tosoh G8 Glycohemoglobin Analyzer,
# Trinity Biotech Premier Hb9210 Automated HPLC System), their generic
# instrument instances, and their deployments.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. InstrumentModels sheet -- add two new model rows (5 and 6)
# ---------------------------------------------------------------
$wsModels = $wb.Worksheets.Item("InstrumentModels")

$wsModels.Range("A5").Value = "nhanes-kb:TOSOH-G8-GLYCOHEMOGLOBIN-ANALYZER"
$wsModels.Range("B5").Value = "vstoi:PhysicalInstrument"
$wsModels.Range("C5").Value = "Tosoh G8 Glycohemoglobin Analyzer"
$wsModels.Range("D5").Value = "Tosoh Bioscience, Inc."

$wsModels.Range("A6").Value = "nhanes-kb:TRINITY-BIOTECH-PREMIER-HB9210-AUTOMATED-HPLC-SYSTEM"
$wsModels.Range("B6").Value = "vstoi:PhysicalInstrument"
$wsModels.Range("C6").Value = "Trinity Biotech Premier Hb9210 Automated HPLC System"
$wsModels.Range("D6").Value = "Trinity Biotech"

$wsModels.Range("A5:D6").Style = "Normal"
$wsModels.Range("C5:C6").Select()

# ---------------------------------------------------------------
# 2. Instruments sheet -- fill in the two generic instrument rows
#    (70 and 71) that already existed as blank placeholder rows.
#    Label column (C) was entered first, then the URI column (A),
#    then the model-reference column (B).
# ---------------------------------------------------------------
$wsInstr = $wb.Worksheets.Item("Instruments")

$wsInstr.Range("C70").Value = "Generic Tosoh G8 Glycohemoglobin Analyzer"
$wsInstr.Range("C71").Value = "Generic Trinity Biotech Premier Hb9210 Automated HPLC System"

$wsInstr.Range("A70").Value = "nhanes-kb:INS-TOSOH-G8-GLYCOHEMOGLOBIN-ANALYZER"
$wsInstr.Range("A71").Value = "nhanes-kb:INS-TRINITY-BIOTECH-PREMIER-HB9210-AUTOMATED-HPLC-SYSTEM"

$wsInstr.Range("B70").Value = "nhanes-kb:TOSOH-G8-GLYCOHEMOGLOBIN-ANALYZER"
$wsInstr.Range("B71").Value = "nhanes-kb:TRINITY-BIOTECH-PREMIER-HB9210-AUTOMATED-HPLC-SYSTEM"

$wsInstr.Range("A70:C71").Select()

# ---------------------------------------------------------------
# 3. Deployments sheet -- add two new deployment rows (25 and 26)
# ---------------------------------------------------------------
$wsDeploy = $wb.Worksheets.Item("Deployments")

$wsDeploy.Range("A25").Value = "nhanes-kb:DPL-TOSOH-G8-GLYCOHEMOGLOBIN-ANALYZER"
$wsDeploy.Range("A26").Value = "nhanes-kb:DPL-TRINITY-BIOTECH-PREMIER-HB9210-AUTOMATED-HPLC-SYSTEM"

$wsDeploy.Range("B25").Value = "vstoi:Deployment"
$wsDeploy.Range("B26").Value = "vstoi:Deployment"

$wsDeploy.Range("C25").Value = "nhanes-kb:PLT-GENERIC-PLATFORM"
$wsDeploy.Range("C26").Value = "nhanes-kb:PLT-GENERIC-PLATFORM"

$wsDeploy.Range("D25").Value = "nhanes-kb:INS-TOSOH-G8-GLYCOHEMOGLOBIN-ANALYZER"
$wsDeploy.Range("D26").Value = "nhanes-kb:INS-TRINITY-BIOTECH-PREMIER-HB9210-AUTOMATED-HPLC-SYSTEM"

$wsDeploy.Range("E25").Value = "nhanes-kb:DET-GENERIC-DETECTOR"
$wsDeploy.Range("E26").Value = "nhanes-kb:DET-GENERIC-DETECTOR"

$wsDeploy.Range("F25").Value = "2015-11-29T11:00:00.999Z"
$wsDeploy.Range("F26").Value = "2015-11-29T11:00:00.999Z"

$wsDeploy.Range("A27").Select()

# ---------------------------------------------------------------
# 4. Autofit the columns whose widest entry changed because of the
#    new, longer strings above.
# ---------------------------------------------------------------
$wsModels.Columns.Item(1).AutoFit()
$wsModels.Columns.Item(3).AutoFit()
$wsInstr.Columns.Item(1).AutoFit()
$wsInstr.Columns.Item(2).AutoFit()
$wsDeploy.Columns.Item(1).AutoFit()
$wsDeploy.Columns.Item(4).AutoFit()

$wsDeploy.Activate()
